{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the two content changes from the diff:\n//   1. In the \"Generating JSON\" section, remove the spurious\n//      w:proofErr gramStart/gramEnd markers that straddled the \"(\"\n//      character, merging \"(\" into the following run so the text\n//      reads \"toJson() method ...\" with a clean run layout.\n//   2. In the \"Photo management stuffs\" section, append a new\n//      sentence after \"...images folder on the device\", split across\n//      two runs exactly as authored:\n//        \".  This will allow us to e\" + \"asily attach the images for\n//        the relevant key locations while being able to just\n//        reference a file path in the JSON string.\"\n//\n// We rebuild each affected paragraph's OOXML explicitly (via\n// Range.insertOoxml with location \"Replace\") so the resulting run /\n// proofErr structure matches the target precisely, rather than\n// relying on insertText (which does not let us drop the proofErr\n// elements or control run boundaries).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"';\n\nfunction pkgOoxml(paragraphXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + W_NS + '><w:body>' + paragraphXml + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>'\n  );\n}\n\n// --- Change 1: \"toJson(\" / \") method ...\" -> drop gramStart/gramEnd ---\nlet jsonParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"All relevant classes will have a\") !== -1) {\n    jsonParagraph = p;\n    break;\n  }\n}\n\nif (jsonParagraph) {\n  // Fall back to the paragraph's known attributes (straight from the\n  // source document) if for some reason they cannot be re-read live.\n  let pOpenTag = '<w:p w14:paraId=\"0094CB7F\" w14:textId=\"511E95C2\" w:rsidR=\"0069762D\" w:rsidRDefault=\"0069762D\" w:rsidP=\"0069762D\">';\n  try {\n    const ooxmlResult = jsonParagraph.getOoxml();\n    await context.sync();\n    const m = ooxmlResult.value.match(/<w:p(?: [^>]*)?>/);\n    if (m) pOpenTag = m[0];\n  } catch (e) {\n    // keep fallback\n  }\n\n  const newParaXml =\n    pOpenTag +\n      '<w:r><w:t xml:space=\"preserve\">All relevant classes will have a </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>toJson</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t>() method that constructs a JSON object based on the c</w:t></w:r>' +\n      '<w:r w:rsidR=\"00E74D22\"><w:t xml:space=\"preserve\">ontents of the class.  This includes any objects constructed by the class itself.  Any photographs that we are going to use have their file path added as a component of the generated JSON objects. </w:t></w:r>' +\n    '</w:p>';\n  jsonParagraph.getRange().insertOoxml(pkgOoxml(newParaXml), \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2: append new sentence after \"...images folder on the device\" ---\nconst body2 = context.document.body;\nconst paragraphs2 = body2.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet photoParagraph = null;\nfor (const p of paragraphs2.items) {\n  if (p.text.indexOf(\"images folder on the device\") !== -1) {\n    photoParagraph = p;\n    break;\n  }\n}\n\nif (photoParagraph) {\n  let pOpenTag = '<w:p w14:paraId=\"0E1E9574\" w14:textId=\"2507E85F\" w:rsidR=\"00AD2AF2\" w:rsidRPr=\"00AD2AF2\" w:rsidRDefault=\"00410526\" w:rsidP=\"00AD2AF2\">';\n  try {\n    const ooxmlResult = photoParagraph.getOoxml();\n    await context.sync();\n    const m = ooxmlResult.value.match(/<w:p(?: [^>]*)?>/);\n    if (m) pOpenTag = m[0];\n  } catch (e) {\n    // keep fallback\n  }\n\n  const newParaXml =\n    pOpenTag +\n      '<w:r><w:t>The android operating system allows us to simply store our own images that the user will create in app in our own file storage system under the images folder on the device</w:t></w:r>' +\n      '<w:r><w:t>.  This will allow us to e</w:t></w:r>' +\n      '<w:r><w:t>asily attach the images for the relevant key locations while being able to just reference a file path in the JSON string.</w:t></w:r>' +\n    '</w:p>';\n  photoParagraph.getRange().insertOoxml(pkgOoxml(newParaXml), \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the two content changes from the diff:\n#   1. In the \"Generating JSON\" section, remove the spurious\n#      w:proofErr gramStart/gramEnd markers that straddled the \"(\"\n#      character, merging \"(\" into the following run so the text\n#      reads \"toJson() method ...\" with a clean run layout.\n#   2. In the \"Photo management stuffs\" section, append a new\n#      sentence after \"...images folder on the device\", split across\n#      two runs exactly as authored:\n#        \".  This will allow us to e\" + \"asily attach the images for\n#        the relevant key locations while being able to just\n#        reference a file path in the JSON string.\"\n#\n# We rebuild each affected paragraph's OOXML explicitly (via\n# Range.InsertXML) so the resulting run / proofErr structure matches\n# the target precisely.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"toJson(\" / \") method ...\" -> drop gramStart/gramEnd ---\n$jsonParaXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"0094CB7F\" w14:textId=\"511E95C2\" w:rsidR=\"0069762D\" w:rsidRDefault=\"0069762D\" w:rsidP=\"0069762D\"><w:r><w:t xml:space=\"preserve\">All relevant classes will have a </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>toJson</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>() method that constructs a JSON object based on the c</w:t></w:r><w:r w:rsidR=\"00E74D22\"><w:t xml:space=\"preserve\">ontents of the class.  This includes any objects constructed by the class itself.  Any photographs that we are going to use have their file path added as a component of the generated JSON objects. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*All relevant classes will have a*\") {\n    $p.Range.InsertXML($jsonParaXml)\n    break\n  }\n}\n\n# --- Change 2: append new sentence after \"...images folder on the device\" ---\n$photoParaXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"0E1E9574\" w14:textId=\"2507E85F\" w:rsidR=\"00AD2AF2\" w:rsidRPr=\"00AD2AF2\" w:rsidRDefault=\"00410526\" w:rsidP=\"00AD2AF2\"><w:r><w:t>The android operating system allows us to simply store our own images that the user will create in app in our own file storage system under the images folder on the device</w:t></w:r><w:r><w:t>.  This will allow us to e</w:t></w:r><w:r><w:t>asily attach the images for the relevant key locations while being able to just reference a file path in the JSON string.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*images folder on the device*\") {\n    $p.Range.InsertXML($photoParaXml)\n    break\n  }\n}\n"}
